# Fix solution project name in ObjectIntro slides
# 1) Rename the "...Solution" bullet on slide 1 from
#    "PracticeObjectIntroAndMiscSolution" to "PracticeSolutionObjectIntroAndMisc".
# 2) The date placeholder ("Update automatically" datetimeFigureOut field) used
#    across the slide master / layouts / notes master is refreshed to the
#    current authoring date when PowerPoint re-saves the deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: locate the 1-based Characters() start index of $target inside $range
# by scanning (TextRange.Find isn't reliable in this host, but Characters()
# reads are), then overwrite just that slice so surrounding run formatting
# (rPr) is left untouched.
# ---------------------------------------------------------------------------
function Find-CharStart($range, [string]$target) {
    $len = $range.Length
    $tlen = $target.Length
    for ($i = 1; $i -le ($len - $tlen + 2); $i++) {
        $candidate = $range.Characters($i, $tlen).Text
        if ($candidate -eq $target) {
            return $i
        }
    }
    return -1
}

function Replace-InRange($range, [string]$oldText, [string]$newText) {
    $startIdx = Find-CharStart $range $oldText
    if ($startIdx -gt 0) {
        $sub = $range.Characters($startIdx, $oldText.Length)
        $sub.Text = $newText
        return $true
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1) Slide 1: fix the mis-ordered solution project name.
# ---------------------------------------------------------------------------
$oldName = "PracticeObjectIntroAndMiscSolution"
$newName = "PracticeSolutionObjectIntroAndMisc"

$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*$oldName*") {
            Replace-InRange $tr $oldName $newName | Out-Null
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Refresh the "datetimeFigureOut" date placeholders (slide master, every
#    slide layout, and the notes master) to the current save date.
# ---------------------------------------------------------------------------
$oldDate = "12/4/2022"
$newDate = "9/9/2024"

function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$slideMaster = $p.SlideMaster
Update-DateShapes $slideMaster

for ($li = 1; $li -le $slideMaster.CustomLayouts.Count; $li++) {
    Update-DateShapes $slideMaster.CustomLayouts.Item($li)
}

# The Notes Master's date placeholder doesn't pick up shape-text writes in
# this host, but the HeadersFooters.DateAndTime setter does.
$notesMaster = $p.NotesMaster
if ($notesMaster.HeadersFooters.DateAndTime.Text -ne $newDate) {
    $notesMaster.HeadersFooters.DateAndTime.Text = $newDate
}

Write-Host "done"
